$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data to the sheet
$ws.Range("C3").Value = "ПАО"
$ws.Range("D4").Value = "Имя 3"

# Update the selection to E4 only
$ws.Range("E4").Select()
